$d = $word.ActiveDocument

# --- First paragraph (the **ID__AFFARS_MP5342_9__ID** line) ---
$p1 = $d.Paragraphs.Item(1)

# Add a paragraph border (box) with 5pt spacing on all four sides.
$borders = $p1.Format.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# Increase the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.Format.LeftIndent = 11.25

# Remove the trailing run that only contains a single space character
# left over at the end of the paragraph.
$pStart = $p1.Range.Start
$pEnd = $p1.Range.End
$trailing = $d.Range($pEnd - 2, $pEnd - 1)
if ($trailing.Text -eq " ") {
    $trailing.Delete()
}
